$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row for the "sports" recommendation endpoint in the previously
# empty row 3 (between the "restaurant" and "music" recommendation rows).
$ws.Range("A3").Value = "/recommendation/sports"
$ws.Range("B3").Value = "get a recommendation of sports, return a list of sports (only 1 now)"
$ws.Rows.Item(3).RowHeight = 15.75

# Fix the misspelling of "restaurant" in the first row.
$ws.Range("A1").Value = "/recommendation/restaurant"

# Reset the selection to the default cell so the saved view doesn't keep the
# cursor pinned at the old C23 position.
$null = $ws.Range("A1").Select()
